$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheets
# ------------------------------------------------------------------
$wsTeam    = $wb.Worksheets.Item("ProjectTeam")
$wsProduct = $wb.Worksheets.Item("Product Backlog")
$wsSprint  = $wb.Worksheets.Item("Sprint Backlog")
$wsBurn    = $wb.Worksheets.Item("BurndownChart")

# ==================================================================
# Product Backlog sheet
# ==================================================================

# Row 5: effort plan original changes from 15 to 20 (text stays the same)
$wsProduct.Range("E5").Value2 = 20

# Rows 6 and 7 swap their entire content (story, description, priority)
$wsProduct.Range("B6").Value2 = "Tagesübersicht einsehen"
$wsProduct.Range("C6").Value2 = "Die tägliche Übersicht der Spitex Mitarbeiter"
$wsProduct.Range("D6").Value2 = "medium"
$wsProduct.Range("E6").Value2 = 15

$wsProduct.Range("B7").Value2 = "Einsatz starten / beenden"
$wsProduct.Range("C7").Value2 = "Den Einsatz eines Spitex Mitarbeiters starten und beenden"
$wsProduct.Range("D7").Value2 = "high"
$wsProduct.Range("E7").Value2 = 30

# Row 10 (Total row with SUM formulas) is removed entirely
$wsProduct.Rows.Item(10).ClearContents()

# ==================================================================
# Sprint Backlog sheet
# ==================================================================

# --- Column widths -------------------------------------------------
$wsSprint.Columns.Item(3).ColumnWidth = 33.1

# --- Row 2 (Domain Model) ------------------------------------------
$wsSprint.Range("C2").Value2 = "Domain Model & Testdaten erstellen"
$wsSprint.Range("F2").Value2 = "Seglias"
$wsSprint.Range("G2").Value2 = "Ritz"
$wsSprint.Range("J2").ClearContents()
$wsSprint.Range("K2").Value2 = 9
$wsSprint.Range("L2").Value2 = "done"

# --- Row 3 (Wochenplanung für MA & Patient) -------------------------
$wsSprint.Range("C3").Value2 = "Wochenplanung für MA & Patient"
$wsSprint.Range("D3").WrapText = $true
$wsSprint.Range("D3").Value2 = "Frontend und Backend für die MA-View und die Patienten-View muss erstellt werden. Es soll ein Kalender dargestellt werden, welcher Termine einer Woche beinhaltet."
$wsSprint.Range("E3").Value2 = "Planung"
$wsSprint.Range("F3").Value2 = "Berger"
$wsSprint.Range("G3").Value2 = "Seglias"
$wsSprint.Range("J3").ClearContents()
$wsSprint.Range("K3").Value2 = 9
$wsSprint.Range("L3").Value2 = "done"
$wsSprint.Rows.Item(3).RowHeight = 57.9

# --- Row 4 (Basis für Layering Pattern legen) -----------------------
$wsSprint.Range("C4").Value2 = "Basis für Layering Pattern legen"
$wsSprint.Range("D4").WrapText = $true
$wsSprint.Range("D4").Value2 = "Business-, Persistenzservice und API sollen vorhanden sein"
$wsSprint.Range("E4").Value2 = "Layering"
$wsSprint.Range("F4").Value2 = "Ritz"
$wsSprint.Range("G4").Value2 = "Berger"
$wsSprint.Range("J4").ClearContents()
$wsSprint.Range("K4").Value2 = 9
$wsSprint.Range("L4").Value2 = "done"
$wsSprint.Rows.Item(4).RowHeight = 28.8

# --- Row 5 (Patienteneinsätze GUI) -----------------------------------
$wsSprint.Range("C5").Value2 = "Patienteneinsätze GUI"
$wsSprint.Range("D5").WrapText = $true
$wsSprint.Range("D5").Value2 = "Die View für das Erfassen von neuen Patienten Einsätzen muss vorhanden sein"
$wsSprint.Range("E5").Value2 = "Planung"
$wsSprint.Range("F5").Value2 = "Nussbaum"
$wsSprint.Range("G5").Value2 = "Schüpbach"
$wsSprint.Range("J5").ClearContents()
$wsSprint.Range("K5").Value2 = 9
$wsSprint.Range("L5").Value2 = "done"
$wsSprint.Rows.Item(5).RowHeight = 28.8

# --- Row 6 (Logik der Patienteneinsätze) ------------------------------
$wsSprint.Range("C6").Value2 = "Logik der Patienteneinsätze"
$wsSprint.Range("D6").WrapText = $true
$wsSprint.Range("D6").Value2 = "Backend Logik für das Erfassen von neuen Einsätzen muss vorhanden sein. Terminkonflikte müssen verhindert werden. "
$wsSprint.Range("E6").Value2 = "Planung"
$wsSprint.Range("F6").Value2 = "Schüpbach"
$wsSprint.Range("G6").Value2 = "Nussbaum"
$wsSprint.Range("J6").ClearContents()
$wsSprint.Range("K6").Value2 = 9
$wsSprint.Range("L6").Value2 = "done"
$wsSprint.Rows.Item(6).RowHeight = 57.6

# --- Row 8 (Frontend für MA Zuweisung) --------------------------------
$wsSprint.Range("A8").Value2 = 3.1
$wsSprint.Range("B8").Value2 = 2
$wsSprint.Range("C8").Value2 = "Frontend für MA Zuweisung"
$wsSprint.Range("D8").WrapText = $true
$wsSprint.Range("D8").Value2 = "MA kann einem Einsatz zugewiesen werden"
$wsSprint.Range("E8").Value2 = "Planung"
$wsSprint.Range("F8").Value2 = "Schüpbach"
$wsSprint.Range("G8").Value2 = "Nussbaum"
$wsSprint.Range("H8").Value2 = "medium"
$wsSprint.Range("I8").Value2 = 5
$wsSprint.Range("K8").Value2 = 0
$wsSprint.Range("L8").Value2 = "waiting"

# --- Row 9 (Backend für MA Zuweisungsvorschläge) -----------------------
$wsSprint.Range("A9").Value2 = 3.2
$wsSprint.Range("B9").Value2 = 2
$wsSprint.Range("C9").Value2 = "Backend für MA Zuweisungsvorschläge"
$wsSprint.Range("D9").WrapText = $true
$wsSprint.Range("D9").Value2 = "Vorschläge nach Priorität sollen vorhanden sein, Terminkonflikte sollen berücksichtigt sein. Die Priorität ist nach Anzahl bisheriger Besuche festgelegt."
$wsSprint.Range("E9").Value2 = "Planung"
$wsSprint.Range("F9").Value2 = "Berger"
$wsSprint.Range("G9").Value2 = "Ritz"
$wsSprint.Range("H9").Value2 = "medium"
$wsSprint.Range("I9").Value2 = 5
$wsSprint.Range("K9").Value2 = 0
$wsSprint.Range("L9").Value2 = "waiting"
$wsSprint.Rows.Item(9).RowHeight = 57.6

# --- Row 10 (Backend für MA Zuweisung) ---------------------------------
$wsSprint.Range("A10").Value2 = 3.3
$wsSprint.Range("B10").Value2 = 2
$wsSprint.Range("C10").Value2 = "Backend für MA Zuweisung"
$wsSprint.Range("D10").WrapText = $true
$wsSprint.Range("D10").Value2 = "Neue Mission Instanzen müssen erstellt werden, sobald einem Termin ein MA zugewiesen wird."
$wsSprint.Range("E10").Value2 = "Planung"
$wsSprint.Range("F10").Value2 = "Berger"
$wsSprint.Range("G10").Value2 = "Ritz"
$wsSprint.Range("H10").Value2 = "medium"
$wsSprint.Range("I10").Value2 = 5
$wsSprint.Range("K10").Value2 = 0
$wsSprint.Range("L10").Value2 = "waiting"
$wsSprint.Rows.Item(10).RowHeight = 43.2

# --- Row 11 (Frontend für Terminvorschläge MA) -------------------------
$wsSprint.Range("A11").Value2 = 4.1
$wsSprint.Range("B11").Value2 = 2
$wsSprint.Range("C11").Value2 = "Frontend für Terminvorschläge MA"
$wsSprint.Range("D11").WrapText = $true
$wsSprint.Range("D11").Value2 = "Auf der MA Übersicht sollen für noch freie Termin Vorschläge gemacht werden. "
$wsSprint.Range("E11").Value2 = "Planung"
$wsSprint.Range("F11").Value2 = "Schüpbach"
$wsSprint.Range("G11").Value2 = "Seglias"
$wsSprint.Range("H11").Value2 = "medium"
$wsSprint.Range("I11").Value2 = 5
$wsSprint.Range("K11").Value2 = 0
$wsSprint.Range("L11").Value2 = "waiting"
$wsSprint.Rows.Item(11).RowHeight = 28.8

# --- Row 12 (Backend für Terminvorschläge MA) ---------------------------
$wsSprint.Range("A12").Value2 = 4.2
$wsSprint.Range("B12").Value2 = 2
$wsSprint.Range("C12").Value2 = "Backend für Terminvorschläge MA"
$wsSprint.Range("D12").WrapText = $true
$wsSprint.Range("D12").Value2 = "Die Vorschläge sollen auf Anzahl bisheriger Besuche und freiem Terminkalender basieren"
$wsSprint.Range("E12").Value2 = "Planung"
$wsSprint.Range("F12").Value2 = "Seglias"
$wsSprint.Range("G12").Value2 = "Schüpbach"
$wsSprint.Range("H12").Value2 = "medium"
$wsSprint.Range("I12").Value2 = 10
$wsSprint.Range("K12").Value2 = 0
$wsSprint.Range("L12").Value2 = "waiting"
$wsSprint.Rows.Item(12).RowHeight = 43.2

# --- Row 13 (Frontend für Terminvorschläge annehmen) --------------------
$wsSprint.Range("A13").Value2 = 4.3
$wsSprint.Range("B13").Value2 = 2
$wsSprint.Range("C13").Value2 = "Frontend für Terminvorschläge annehmen"
$wsSprint.Range("D13").WrapText = $true
$wsSprint.Range("D13").Value2 = "Die Vorschläge sollen akzeptiert werden können. Die Akzeptierung des Vorschlags generiert eine neue Mission (Wie in 3.3)"
$wsSprint.Range("E13").Value2 = "Planung"
$wsSprint.Range("F13").Value2 = "Nussbaum"
$wsSprint.Range("G13").Value2 = "Schüpbach"
$wsSprint.Range("H13").Value2 = "medium"
$wsSprint.Range("I13").Value2 = 5
$wsSprint.Range("K13").Value2 = 0
$wsSprint.Range("L13").Value2 = "waiting"
$wsSprint.Rows.Item(13).RowHeight = 43.2

# --- Row 14 (Frontend für Tagesübersicht MA erstellen) -------------------
$wsSprint.Range("A14").Value2 = 5.1
$wsSprint.Range("B14").Value2 = 2
$wsSprint.Range("C14").Value2 = "Frontend für Tagesübersicht MA erstellen"
$wsSprint.Range("D14").WrapText = $true
$wsSprint.Range("D14").Value2 = "Alle Termine des aktuellen Tages sollen dargestellt werden. (Timeline)               Einsatz starten & Einsatz beenden sollen vorhanden sein (nicht funktionstüchtig)"
$wsSprint.Range("E14").Value2 = "HealthVisitor"
$wsSprint.Range("F14").Value2 = "Nussbaum"
$wsSprint.Range("G14").Value2 = "Berger"
$wsSprint.Range("H14").Value2 = "medium"
$wsSprint.Range("I14").Value2 = 10
$wsSprint.Range("K14").Value2 = 0
$wsSprint.Range("L14").Value2 = "waiting"
$wsSprint.Rows.Item(14).RowHeight = 57.6

# --- Row 15 (Backend für Tagesübersicht eines MA) -------------------------
$wsSprint.Range("A15").Value2 = 5.2
$wsSprint.Range("B15").Value2 = 2
$wsSprint.Range("C15").Value2 = "Backend für Tagesübersicht eines MA"
$wsSprint.Range("D15").WrapText = $true
$wsSprint.Range("D15").Value2 = "Die Daten die auf der View dargestellt werden, müssen geliefert werden."
$wsSprint.Range("E15").Value2 = "HealthVisitor"
$wsSprint.Range("F15").Value2 = "Ritz"
$wsSprint.Range("G15").Value2 = "Berger"
$wsSprint.Range("H15").Value2 = "medium"
$wsSprint.Range("I15").Value2 = 5
$wsSprint.Range("K15").Value2 = 0
$wsSprint.Range("L15").Value2 = "waiting"
$wsSprint.Rows.Item(15).RowHeight = 28.8

# ==================================================================
# View state: active cell per sheet + active tab
# ==================================================================
$wsTeam.Activate()

$wsProduct.Activate()
$wsProduct.Range("C3").Select() | Out-Null

$wsBurn.Activate()
$wsBurn.Range("E21").Select() | Out-Null

$wsSprint.Activate()
$wsSprint.Range("L8").Select() | Out-Null
